$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.861.52'
$ws.Range("E2").Value = '  -4.38%  '
$ws.Range("D3").Value = '2.330.89'
$ws.Range("E3").Value = '  -5.84%  '
$ws.Range("D5").Value = '307.22'
$ws.Range("E5").Value = '  -4.11%  '
$ws.Range("D6").Value = '83.87'
$ws.Range("E6").Value = '  -9.14%  '
$ws.Range("D7").Value = '0.528'
$ws.Range("E7").Value = '  -4.12%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '0.482'
$ws.Range("E9").Value = '  -5.82%  '
$ws.Range("D10").Value = '0.0816'
$ws.Range("E10").Value = '  -5.39%  '
$ws.Range("D11").Value = '29.84'
$ws.Range("E11").Value = '  -9.73%  '
$ws.Range("E12").Value = '  -0.82%  '
$ws.Range("D13").Value = '2.703.87'
$ws.Range("E13").Value = '  -5.38%  '
$ws.Range("D14").Value = '6.38'
$ws.Range("E14").Value = '  -7.47%  '
$ws.Range("D15").Value = '14.77'
$ws.Range("E15").Value = '  -4.82%  '
$ws.Range("D16").Value = '2.361.10'
$ws.Range("E16").Value = '  -4.78%  '
$ws.Range("D17").Value = '0.745'
$ws.Range("E17").Value = '  -6.21%  '
$ws.Range("D18").Value = '39.897.89'
$ws.Range("E18").Value = '  -4.12%  '
$ws.Range("D19").Value = '0.0₃0898'
$ws.Range("E19").Value = '  -4.71%  '
$ws.Range("D20").Value = '6.04'
$ws.Range("E20").Value = '  -6.28%  '
$ws.Range("D21").Value = '67.61'
$ws.Range("E21").Value = '  -4.34%  '
$ws.Range("D22").Value = '10.53'
$ws.Range("E22").Value = '  -6.21%  '
$ws.Range("D23").Value = '233.64'
$ws.Range("E23").Value = '  -2.47%  '
$ws.Range("D24").Value = '2.53'
$ws.Range("E24").Value = '  -7.93%  '
$ws.Range("E25").Value = '  -0.03%  '
$ws.Range("E26").Value = '  -8.20%  '
$ws.Range("D27").Value = '23.27'
$ws.Range("E27").Value = '  -7.04%  '
$ws.Range("D28").Value = '2.20'
$ws.Range("E28").Value = '  -2.22%  '
$ws.Range("D29").Value = '9.15'
$ws.Range("E29").Value = '  -6.18%  '
$ws.Range("D30").Value = '34.08'
$ws.Range("E30").Value = '  -6.88%  '
$ws.Range("D31").Value = '150.89'
$ws.Range("E31").Value = '  -4.03%  '
$ws.Range("E32").Value = '  +0.08%  '
$ws.Range("D33").Value = '5.10'
$ws.Range("E33").Value = '  -6.09%  '
$ws.Range("D34").Value = '0.0721'
$ws.Range("E34").Value = '  -5.56%  '
$ws.Range("D35").Value = '2.41'
$ws.Range("E35").Value = '  -5.65%  '
$ws.Range("D36").Value = '0.113'
$ws.Range("E36").Value = '  -2.84%  '
$ws.Range("D37").Value = '2.75'
$ws.Range("E37").Value = '  -4.88%  '
$ws.Range("D38").Value = '0.0988'
$ws.Range("E38").Value = '  -4.89%  '
$ws.Range("D39").Value = '15.60'
$ws.Range("E39").Value = '  -8.91%  '
$ws.Range("D40").Value = '1.68'
$ws.Range("E40").Value = '  -8.61%  '
$ws.Range("D41").Value = '3.78'
$ws.Range("E41").Value = '  -5.92%  '
$ws.Range("D42").Value = '2.35'
$ws.Range("E42").Value = '  -4.42%  '
$ws.Range("D43").Value = '1.944.21'
$ws.Range("E43").Value = '  -2.74%  '
$ws.Range("B44").Value = 'VeChain'
$ws.Range("C44").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D44").Value = '0.0264'
$ws.Range("E44").Value = '  -7.24%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '17.48'
$ws.Range("E45").Value = '  -6.35%  '
$ws.Range("D46").Value = '9.34'
$ws.Range("E46").Value = '  -1.86%  '
$ws.Range("D47").Value = '2.63'
$ws.Range("E47").Value = '  -11.63%  '
$ws.Range("D48").Value = '2.582.44'
$ws.Range("E48").Value = '  -6.09%  '
$ws.Range("D49").Value = '91.84'
$ws.Range("E49").Value = '  -5.84%  '
$ws.Range("D50").Value = '70.57'
$ws.Range("E50").Value = '  -6.87%  '
$ws.Range("D51").Value = '63.29'
$ws.Range("E51").Value = '  -6.06%  '
